$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update "last updated" timestamp text (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 26 de Mayo de 2020 a las 08:05"

# --- Swap country names whose order changed in the shared-strings table ---
# Hungria now sorts before Guatemala
$ws.Range("A74").Value = "Hungria"
$ws.Range("A75").Value = "Guatemala"

# El Salvador now sorts before Republica de Macedonia
$ws.Range("A89").Value = "El Salvador"
$ws.Range("A90").Value = "Republica de Macedonia"

# --- Update statistic values (refreshed data) ---
$ws.Range("D11").Value = 162000
$ws.Range("E11").Value = 10361

$ws.Range("B68").Value = 5044
$ws.Range("C68").Value = 154
$ws.Range("D68").Value = 1917
$ws.Range("E68").Value = 2956
$ws.Range("G68").Value = 6
$ws.Range("H68").Value = 171

$ws.Range("B74").Value = 3771
$ws.Range("C74").Value = 15
$ws.Range("D74").Value = 1836
$ws.Range("E74").Value = 1436
$ws.Range("G74").Value = 8
$ws.Range("H74").Value = 499

$ws.Range("B75").Value = 3760
$ws.Range("C75").Value = 336
$ws.Range("D75").Value = 274
$ws.Range("E75").Value = 3427
$ws.Range("G75").Value = 1
$ws.Range("H75").Value = 59

$ws.Range("B77").Value = 3261
$ws.Range("C77").Value = 72
$ws.Range("E77").Value = 641

$ws.Range("B80").Value = 3045
$ws.Range("C80").Value = 3
$ws.Range("D80").Value = 2929
$ws.Range("E80").Value = 59

$ws.Range("B83").Value = 2443
$ws.Range("C83").Value = 10
$ws.Range("D83").Value = 880
$ws.Range("E83").Value = 1433

$ws.Range("B89").Value = 2042
$ws.Range("C89").Value = 59
$ws.Range("D89").Value = 787
$ws.Range("E89").Value = 1219
$ws.Range("G89").Value = 1
$ws.Range("H89").Value = 36

$ws.Range("B90").Value = 1999
$ws.Range("D90").Value = 1439
$ws.Range("E90").Value = 447
$ws.Range("H90").Value = 113

$ws.Range("D162").Value = 37
$ws.Range("E162").Value = 104
